$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 23.76874
$ws.Range("H2").Value = 71.30622000000001
$ws.Range("I2").Value = 0.9555483638834226
$ws.Range("J2").Value = 0.9555483638834227
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 25.15544366666667
$ws.Range("N2").Value = 75.466331
$ws.Range("O2").Value = 0.9701024243751556
$ws.Range("P2").Value = 0.9701024243751556
$ws.Range("Q2").Value = 597.9132000976468
$ws.Range("R2").Value = 5381.218800878821
$ws.Range("S2").Value = 0.9269797844110216
$ws.Range("T2").Value = 0.9269797844110217
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 23.76874
$ws.Range("H3").Value = 71.30622000000001
$ws.Range("I3").Value = 0.9555483638834226
$ws.Range("J3").Value = 0.9555483638834227
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.6301496666666666
$ws.Range("N3").Value = 1.890449
$ws.Range("O3").Value = 0.02430128951224074
$ws.Range("P3").Value = 0.02430128951224074
$ws.Range("Q3").Value = 14.97786358808667
$ws.Range("R3").Value = 134.80077229278
$ws.Range("S3").Value = 0.02322105743367902
$ws.Range("T3").Value = 0.02322105743367902
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 23.76874
$ws.Range("H4").Value = 71.30622000000001
$ws.Range("I4").Value = 0.9555483638834226
$ws.Range("J4").Value = 0.9555483638834227
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1451156666666667
$ws.Range("N4").Value = 0.435347
$ws.Range("O4").Value = 0.005596286112603657
$ws.Range("P4").Value = 0.005596286112603657
$ws.Range("Q4").Value = 3.449216550926667
$ws.Range("R4").Value = 31.04294895834
$ws.Range("S4").Value = 0.005347522038721944
$ws.Range("T4").Value = 0.005347522038721944
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.7500946666666666
$ws.Range("H5").Value = 2.250284
$ws.Range("I5").Value = 0.03015522621270687
$ws.Range("J5").Value = 0.03015522621270688
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 25.15544366666667
$ws.Range("N5").Value = 75.466331
$ws.Range("O5").Value = 0.9701024243751556
$ws.Range("P5").Value = 0.9701024243751556
$ws.Range("Q5").Value = 18.86896413200044
$ws.Range("R5").Value = 169.820677188004
$ws.Range("S5").Value = 0.02925365805652818
$ws.Range("T5").Value = 0.02925365805652818
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.7500946666666666
$ws.Range("H6").Value = 2.250284
$ws.Range("I6").Value = 0.03015522621270687
$ws.Range("J6").Value = 0.03015522621270688
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.6301496666666666
$ws.Range("N6").Value = 1.890449
$ws.Range("O6").Value = 0.02430128951224074
$ws.Range("P6").Value = 0.02430128951224074
$ws.Range("Q6").Value = 0.4726719041684443
$ws.Range("R6").Value = 4.254047137515999
$ws.Range("S6").Value = 0.0007328108825021007
$ws.Range("T6").Value = 0.0007328108825021008
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.7500946666666666
$ws.Range("H7").Value = 2.250284
$ws.Range("I7").Value = 0.03015522621270687
$ws.Range("J7").Value = 0.03015522621270688
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.1451156666666667
$ws.Range("N7").Value = 0.435347
$ws.Range("O7").Value = 0.005596286112603657
$ws.Range("P7").Value = 0.005596286112603657
$ws.Range("Q7").Value = 0.1088504876164444
$ws.Range("R7").Value = 0.9796543885479998
$ws.Range("S7").Value = 0.0001687572736765933
$ws.Range("T7").Value = 0.0001687572736765933
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.3556153333333333
$ws.Range("H8").Value = 1.066846
$ws.Range("I8").Value = 0.01429640990387057
$ws.Range("J8").Value = 0.01429640990387057
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 25.15544366666667
$ws.Range("N8").Value = 75.466331
$ws.Range("O8").Value = 0.9701024243751556
$ws.Range("P8").Value = 0.9701024243751556
$ws.Range("Q8").Value = 8.945661484669555
$ws.Range("R8").Value = 80.51095336202599
$ws.Range("S8").Value = 0.01386898190760582
$ws.Range("T8").Value = 0.01386898190760582
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.3556153333333333
$ws.Range("H9").Value = 1.066846
$ws.Range("I9").Value = 0.01429640990387057
$ws.Range("J9").Value = 0.01429640990387057
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.6301496666666666
$ws.Range("N9").Value = 1.890449
$ws.Range("O9").Value = 0.02430128951224074
$ws.Range("P9").Value = 0.02430128951224074
$ws.Range("Q9").Value = 0.2240908837615555
$ws.Range("R9").Value = 2.016817953854
$ws.Range("S9").Value = 0.0003474211960596246
$ws.Range("T9").Value = 0.0003474211960596246
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.3556153333333333
$ws.Range("H10").Value = 1.066846
$ws.Range("I10").Value = 0.01429640990387057
$ws.Range("J10").Value = 0.01429640990387057
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.1451156666666667
$ws.Range("N10").Value = 0.435347
$ws.Range("O10").Value = 0.005596286112603657
$ws.Range("P10").Value = 0.005596286112603657
$ws.Range("Q10").Value = 0.05160535617355556
$ws.Range("R10").Value = 0.4644482055619999
$ws.Range("S10").Value = 0.00008000680020512026
$ws.Range("T10").Value = 0.00008000680020512026
